# Fixed update to excel issue
# - Rename the "Requested quantity" header on the two existing sheets to
#   more descriptive, code-friendly names.
# - Add a new "PO Forecast" worksheet with the forecast series (ds,
#   PO_Forecast, yhat_lower, yhat_upper), reusing the same header/date
#   styling already used on the other sheets.

$wb = $excel.ActiveWorkbook

# --- Weekly Quantity sheet: rename header B1 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "Weekly_PO_Qty"

# --- Monthly Trend sheet: rename header B1 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "PO Forecast"

# Reuse the bold/bordered/centered header style from the existing sheets
# for the new header row.
$ws1.Range("B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

$forecastRows = @(
    @(45585.99999999999, 40, 39.99999994756828, 40.00000005078284),
    @(45592.99999999999, 40, 39.99999994580869, 40.00000005014621),
    @(45627.99999999999, 40, 39.99999994818493, 40.00000004856577),
    @(45634.99999999999, 40, 39.99999994507549, 40.0000000569318),
    @(45641.99999999999, 40, 39.99999992809531, 40.00000006902613),
    @(45648.99999999999, 40, 39.99999990989347, 40.00000008572911),
    @(45655.99999999999, 40, 39.99999986248793, 40.00000013074359),
    @(45662.99999999999, 40, 39.99999976437518, 40.00000019356222),
    @(45669.99999999999, 40, 39.99999968636185, 40.00000027427507),
    @(45676.99999999999, 40, 39.99999959336695, 40.00000039674766),
    @(45683.99999999999, 40, 39.99999942827622, 40.00000049835391)
)

# Reuse the date-formatted style (column A on the other sheets) for the
# "ds" column of the new sheet.
$ws1.Range("A2").Copy()
$ws3.Range("A2:A12").PasteSpecial(-4122)

$r = 2
foreach ($row in $forecastRows) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Keep the first sheet active/selected, matching the original workbook view.
$ws1.Activate()
$ws1.Range("A1").Select()
